# Fix: prevent hidden columns (here: the "Änderung" / change-marker column L)
# from being labeled when a change is detected, and apply the correct
# "Segmentgruppe" header-row styling to the rows that start a new SG group.
#
# Behaviour implemented:
#  1. Rows that start a new Segmentgruppe-group (13,17,23,27,34,40,80) get the
#     same grey "group header" fill/bold formatting that already exists on
#     rows 2 and 9 (template rows) across the whole A:V range.
#  2. Column L (the hidden "Änderung" column) gets its "ÄNDERUNG" label and
#     highlight style cleared for every row of the affected groups
#     (rows 13-43 and row 80), since that column is hidden and must not be
#     labeled.
#
# Formatting is applied by copying (format-only) from already-correct
# template cells (row 2 / L2) so no new style/font records are introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that begin a new Segmentgruppe and need the full-row "group header" style
$headerRows = @(13,17,23,27,34,40,80)

# All rows in the affected groups whose column L must lose the ÄNDERUNG label
# (this includes both the group-header rows and the regular in-group rows)
$allRows = @(13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,80)

# --- Step 1: apply the "group header" formatting (A:V) to header rows -------
# Row 2 is an already-correctly-styled group-header row: column B is bold
# (style used for the segment-name cell), every other column in A:V uses the
# plain grey group-header style.
$headerTarget = $null
foreach ($r in $headerRows) {
    $rowRange = $ws.Range("A" + $r + ":V" + $r)
    if ($headerTarget -eq $null) {
        $headerTarget = $rowRange
    } else {
        $headerTarget = $excel.Union($headerTarget, $rowRange)
    }
}

$ws.Range("A2:V2").Copy()
foreach ($area in $headerTarget.Areas) {
    $area.PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = $false

# --- Step 2: clear column L (hidden "Änderung" column) for all affected rows
$lTarget = $null
foreach ($r in $allRows) {
    $cell = $ws.Range("L" + $r)
    if ($lTarget -eq $null) {
        $lTarget = $cell
    } else {
        $lTarget = $excel.Union($lTarget, $cell)
    }
}

# L2 already carries the correct "cleared" style for this column
$ws.Range("L2").Copy()
foreach ($area in $lTarget.Areas) {
    $area.PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = $false

foreach ($r in $allRows) {
    $ws.Range("L" + $r).Value = ""
}
